# B1-and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table ("Type of document / Definition / ...") is switched from
#    the deck's custom table style to the built-in "Medium Style 2 - Accent 1"
#    gallery style.
# 2) The presentation's colour theme is switched from the "Integral" (Red
#    Violet) theme to the standard "Office Theme" (Office) colour palette -
#    i.e. what you get by picking a different theme on the Design tab. That
#    is exposed on the object model through ThemeColorScheme.Colors(n).RGB,
#    which is the documented automation surface for editing a theme's colour
#    scheme (PowerPoint has no scriptable way to rename the underlying theme
#    / colour-scheme, so only the twelve colour slots are updated here).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 5
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FAE60D27-4496-4860-8CB5-DA5B90442BA2}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour palette -> Office Theme
# ---------------------------------------------------------------------------
# RGB() builds the COM colour (stored as 0xBBGGRR) from a standard RRGGBB hex
# string, the same way PowerPoint's own RGB(r,g,b) helper does.
function HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order of ThemeColorScheme.Colors(): dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToComRGB($officeThemeColors[$i - 1])
}
